# Commit: "updated the file with my details"
# Adds a second row to Sheet1 with the contributor's name, email and repo
# link, turning the email/link cells into real hyperlinks (Excel's
# built-in "Hyperlink" cell style - underlined, theme color 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the literal cell text first, then attach the hyperlinks without an
# explicit TextToDisplay so Excel keeps the text exactly as typed (no
# "mailto:" prefix bleeding into the displayed/stored string).
$ws.Range("A2").Value = "karim abdelmnem mohamed"
$ws.Range("B2").Value = "karimabdelmnem71@gmail.com"
$ws.Range("C2").Value = "https://github.com/Karim3bdelmn3m/our-project.git"

[void]$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:karimabdelmnem71@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/Karim3bdelmn3m/our-project.git")

# Leave the selection where the author left it after entering the data.
[void]$ws.Range("D2").Select()
